$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.09138
$ws.Range("M2").Value = 12.40685866666667
$ws.Range("N2").Value = 37.220576
$ws.Range("O2").Value = 0.1720325859617629
$ws.Range("P2").Value = 0.1720325859617629
$ws.Range("Q2").Value = 0.3779129149866667
$ws.Range("R2").Value = 3.40121623488
$ws.Range("S2").Value = 0.1720325859617629
$ws.Range("T2").Value = 0.1720325859617629

# Row 3
$ws.Range("H3").Value = 0.09138
$ws.Range("O3").Value = 0.6097142007069145
$ws.Range("P3").Value = 0.6097142007069145
$ws.Range("S3").Value = 0.6097142007069145
$ws.Range("T3").Value = 0.6097142007069145

# Row 4
$ws.Range("H4").Value = 0.09138
$ws.Range("N4").Value = 47.220765
$ws.Range("O4").Value = 0.2182532133313226
$ws.Range("P4").Value = 0.2182532133313226
$ws.Range("Q4").Value = 0.4794481673
$ws.Range("S4").Value = 0.2182532133313226
$ws.Range("T4").Value = 0.2182532133313226
